$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the AIP-beats-HTB rows (65-103) as "yes" (was "yes/no") ---
$ws.Range("B65:B103").Value = "yes"

# --- Fill in feature/java file reference pairs per hand-rank group ---
# (written in row order so new shared-string entries come out in the same
# order as the authored workbook)

# AIP has royal flush holds and beats HTB ... (rows 65-73)
$ws.Range("E65:E73").Value = "AIP_royal_flush_beats_HTB.feature"
$ws.Range("F65:F73").Value = "AIPRoyalFlushBeatsHTB.java"

# AIP has straight flush holds and beats HTB ... (rows 74-81)
$ws.Range("E74:E81").Value = "AIP_straight_flush_beats_HTB.feature"
$ws.Range("F74:F81").Value = "AIPStraightFlushBeatsHTB.java"

# AIP has 4-of-a-kind holds and beats HTB ... (rows 82-88)
$ws.Range("E82:E88").Value = "AIP_four_of_a_kind_beats_HTB.feature"
$ws.Range("F82:F88").Value = "AIPFourOfAKindBeatsHTB.java"

# AIP has full house holds and beats HTB ... (rows 89-94)
$ws.Range("E89:E94").Value = "AIP_full_house_beats_HTB.feature"
$ws.Range("F89:F94").Value = "AIPFullHouseBeatsHTB.java"

# AIP has flush holds and beats HTB ... (rows 95-99)
$ws.Range("E95:E99").Value = "AIP_flush_beats_HTB.feature"
$ws.Range("F95:F99").Value = "AIPFlushBeatsHTB.java"

# AIP has straight holds and beats HTB ... (rows 100-103)
$ws.Range("E100:E103").Value = "AIP_straight_beats_HTB.feature"
$ws.Range("F100:F103").Value = "AIPStraightBeatsHTB.java"

# --- New column headers for the "AIP holds and beats HTB" block ---
$ws.Range("E7").Value = "feature file"
$ws.Range("F7").Value = "java file"

# --- Update the view: scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 100
$ws.Range("E7").Select()
